# Weekly fruit/vegetable price update: a new week's data row is inserted
# at row 16 (pushing the existing rows 16-71 down to 17-72), and the new
# row 16 is populated with this week's observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 16, shifting rows 16..71 down to 17..72.
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new weekly record.
$ws.Cells.Item(16, 1).Value2  = 7
$ws.Cells.Item(16, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(16, 3).Value2  = "Ñuble"
$ws.Cells.Item(16, 4).Value2  = 44910
$ws.Cells.Item(16, 5).Value2  = 16
$ws.Cells.Item(16, 6).Value2  = 100112026
$ws.Cells.Item(16, 7).Value2  = "Haba"
$ws.Cells.Item(16, 8).Value2  = "Sin especificar"
$ws.Cells.Item(16, 9).Value2  = "Primera"
$ws.Cells.Item(16, 10).Value2 = 100
$ws.Cells.Item(16, 11).Value2 = 10000
$ws.Cells.Item(16, 12).Value2 = 12000
$ws.Cells.Item(16, 13).Value2 = 11000
$ws.Cells.Item(16, 14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(16, 15).Value2 = "Provincia de Diguillín"
$ws.Cells.Item(16, 16).Value2 = 440
$ws.Cells.Item(16, 17).Value2 = 25
$ws.Cells.Item(16, 18).Value2 = "Hortaliza"
